$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = "Yacht Rental"
$ws.Range("F4").Value = 8000
$ws.Range("F4").NumberFormat = $ws.Range("D5").NumberFormat

$ws.Range("E9").Select()
